# Add a new "ridge_classifier" model section to the scores sheet, mirroring
# the layout of the existing "multinomal_nb" section (rows 47-68).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone the formatting (fonts/fills/borders/number formats/merges aren't
#    carried by PasteSpecial formats, so merges are re-applied below) of the
#    last model block (rows 47:68) onto the new block (rows 70:91).
$ws.Range("A47:F68").Copy()
$ws.Range("A70").PasteSpecial(-4122)
$ws.Range("A47").Copy()

# 2) Re-create the merged ranges for the new block.
$ws.Range("A70:E70").Merge()
$ws.Range("A71:E71").Merge()
$ws.Range("C80:E80").Merge()
$ws.Range("A81:E81").Merge()
$ws.Range("A82:E82").Merge()
$ws.Range("C91:E91").Merge()

# 3) Fix the one formatting inconsistency carried over from the template:
#    B80 ("macro avg f1-score" value) uses the plain numeric style (like
#    B68/B91), not the bold-row numeric style (B57) it inherited by being
#    copied from B57 in step 1.
$ws.Range("B68").Copy()
$ws.Range("B80").PasteSpecial(-4122)

# 4) Write the header / time cells.
$ws.Cells.Item(70, 1).Value = "ridge_classifier (upsamlping + stratified_split + max_tfidf_features)"
$ws.Cells.Item(70, 6).Value = "~20min"
$ws.Cells.Item(71, 1).Value = "appearance"

# 5) First classification-report table (rows 72-80): 0/1/2 classes + accuracy
#    + macro avg + weighted avg + macro avg f1-score summary.
$ws.Cells.Item(72, 1).Value = "label"
$ws.Cells.Item(72, 2).Value = "precision"
$ws.Cells.Item(72, 3).Value = "recall"
$ws.Cells.Item(72, 4).Value = "f1-score"
$ws.Cells.Item(72, 5).Value = "support"

$ws.Cells.Item(73, 1).Value = "0 (Negative)"
$ws.Cells.Item(73, 2).Value = 0.99
$ws.Cells.Item(73, 3).Value = 1
$ws.Cells.Item(73, 4).Value = 0.99
$ws.Cells.Item(73, 5).Value = 166837

$ws.Cells.Item(74, 1).Value = "1 (Neutral)"
$ws.Cells.Item(74, 2).Value = 0.83
$ws.Cells.Item(74, 3).Value = 0.83
$ws.Cells.Item(74, 4).Value = 0.83
$ws.Cells.Item(74, 5).Value = 166838

$ws.Cells.Item(75, 1).Value = "2 (Positive)"
$ws.Cells.Item(75, 2).Value = 0.84
$ws.Cells.Item(75, 3).Value = 0.83
$ws.Cells.Item(75, 4).Value = 0.83
$ws.Cells.Item(75, 5).Value = 166837

# row 76 stays blank (separator row, formatting only)

$ws.Cells.Item(77, 1).Value = "accuracy"
$ws.Cells.Item(77, 4).Value = 0.89
$ws.Cells.Item(77, 5).Value = 500512

$ws.Cells.Item(78, 1).Value = "macro avg"
$ws.Cells.Item(78, 2).Value = 0.88
$ws.Cells.Item(78, 3).Value = 0.89
$ws.Cells.Item(78, 4).Value = 0.89
$ws.Cells.Item(78, 5).Value = 500512

$ws.Cells.Item(79, 1).Value = "weighted avg"
$ws.Cells.Item(79, 2).Value = 0.88
$ws.Cells.Item(79, 3).Value = 0.89
$ws.Cells.Item(79, 4).Value = 0.89
$ws.Cells.Item(79, 5).Value = 500512

$ws.Cells.Item(80, 1).Value = "macro avg f1-score"
$ws.Cells.Item(80, 2).Value = 0.88480000000000003

# row 81 stays blank (separator row, formatting only)

# 6) Second classification-report table ("palate", rows 82-91).
$ws.Cells.Item(82, 1).Value = "palate"

$ws.Cells.Item(83, 1).Value = "label"
$ws.Cells.Item(83, 2).Value = "precision"
$ws.Cells.Item(83, 3).Value = "recall"
$ws.Cells.Item(83, 4).Value = "f1-score"
$ws.Cells.Item(83, 5).Value = "support"

$ws.Cells.Item(84, 1).Value = "0 (Negative)"
$ws.Cells.Item(84, 2).Value = 0.96
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(84, 4).Value = 0.98
$ws.Cells.Item(84, 5).Value = 166838

$ws.Cells.Item(85, 1).Value = "1 (Neutral)"
$ws.Cells.Item(85, 2).Value = 0.85
$ws.Cells.Item(85, 3).Value = 0.81
$ws.Cells.Item(85, 4).Value = 0.83
$ws.Cells.Item(85, 5).Value = 166837

$ws.Cells.Item(86, 1).Value = "2 (Positive)"
$ws.Cells.Item(86, 2).Value = 0.84
$ws.Cells.Item(86, 3).Value = 0.85
$ws.Cells.Item(86, 4).Value = 0.85
$ws.Cells.Item(86, 5).Value = 166837

# row 87 stays blank (separator row, formatting only)

$ws.Cells.Item(88, 1).Value = "accuracy"
$ws.Cells.Item(88, 4).Value = 0.89
$ws.Cells.Item(88, 5).Value = 500512

$ws.Cells.Item(89, 1).Value = "macro avg"
$ws.Cells.Item(89, 2).Value = 0.88
$ws.Cells.Item(89, 3).Value = 0.89
$ws.Cells.Item(89, 4).Value = 0.88
$ws.Cells.Item(89, 5).Value = 500512

$ws.Cells.Item(90, 1).Value = "weighted avg"
$ws.Cells.Item(90, 2).Value = 0.88
$ws.Cells.Item(90, 3).Value = 0.89
$ws.Cells.Item(90, 4).Value = 0.88
$ws.Cells.Item(90, 5).Value = 500512

$ws.Cells.Item(91, 1).Value = "macro avg f1-score"
$ws.Cells.Item(91, 2).Value = 0.88449999999999995

# 7) Point the view back at the new content, matching the author's final
#    selection.
$ws.Range("B95").Select()
